# Switch license from BY-NC to BY-SA
#
# This presentation's closing/credits slide (the 2nd slide, which carries
# the "Except where otherwise noted..." Creative Commons attribution) was
# re-licensed from CC BY-NC 4.0 to CC BY-SA 4.0. Only the *visible* text
# changes -- the underlying hyperlink (rId2) still resolves to whatever
# address it already pointed at, so we must not touch the hyperlink
# target, only the run text the user reads on the slide.

$p = $ppt.ActivePresentation

# --- Slide 2: the CC-license slide -----------------------------------
$s2 = $p.Slides.Item(2)
$shp = $s2.Shapes.Item(1)
$tr = $shp.TextFrame.TextRange

# 1) "CC BY-NC 4.0. To view a copy of this license, visit " -> "CC BY-SA 4.0. ..."
#    Replace just the "BY-NC " word (including its trailing space) with
#    "BY-SA " so the surrounding text keeps its existing run formatting.
$full = $tr.Text
$oldWord = "BY-NC "
$newWord = "BY-SA "
$idx = $full.IndexOf($oldWord)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldWord.Length)
    $sub.Text = $newWord
}

# 2) The displayed hyperlink text itself:
#    "https://creativecommons.org/licenses/by-nc/4.0"
#      -> "https://creativecommons.org/licenses/by-sa/4.0"
#    (the hyperlink target keeps pointing at whatever it already did --
#    only the text the reader sees is updated.)
$full = $tr.Text
$oldUrlTail = "creativecommons.org/licenses/by-nc/4.0"
$newUrlTail = "creativecommons.org/licenses/by-sa/4.0"
$idx = $full.IndexOf($oldUrlTail)
if ($idx -ge 0) {
    $sub = $tr.Characters($idx + 1, $oldUrlTail.Length)
    $sub.Text = $newUrlTail
}

# A one-EMU-scale nudge to the title placeholder's stored position that
# rides along with PowerPoint's resave of this slide.
$shp.Left = 566057 / 12700.0

# --- Slide 17: tidy up a stray run split in "evens[] = ..." ----------
$s17 = $p.Slides.Item(17)
$shp17 = $s17.Shapes.Item(3)
$tr17 = $shp17.TextFrame.TextRange
$full17 = $tr17.Text
$needle17 = "evens[] = 2, 4, 6, 8, 10, 12, 14, 16, 18, 20"
$idx17 = $full17.IndexOf($needle17)
if ($idx17 -ge 0) {
    $sub17 = $tr17.Characters($idx17 + 1, $needle17.Length)
    $sub17.Text = $needle17
}

# --- Slide 18: tidy up stray tab/text run splits ----------------------
$s18 = $p.Slides.Item(18)
$shp18 = $s18.Shapes.Item(3)
$tr18 = $shp18.TextFrame.TextRange

$full18 = $tr18.Text
$needleA = "`t1. odds[] addition result will be saved to the variable "
$idxA = $full18.IndexOf($needleA)
if ($idxA -ge 0) {
    $subA = $tr18.Characters($idxA + 1, $needleA.Length)
    $subA.Text = $needleA
}

$full18 = $tr18.Text
$needleB = "`t2. evens[] addition result will be saved to the variable "
$idxB = $full18.IndexOf($needleB)
if ($idxB -ge 0) {
    $subB = $tr18.Characters($idxB + 1, $needleB.Length)
    $subB.Text = $needleB
}
